# Gestione: pulizia e aggiornamento note, appunti etc
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7: another "Ink - Generale" / "Programmazione" entry
$ws.Range("A7").Value = 45647
$ws.Range("B7").Value = "Programmazione"
$ws.Range("C7").Value = "Ink - Generale"
$ws.Range("D7").Value = 0.125
$ws.Range("E7").Value = 'Ora la struttura ragiona partendo da "main"'

# New row 8: extra note line
$ws.Range("E8").Value = "Sistemati piccoli bug e simili"

# Move the active selection to E10 (matches the saved view state)
$ws.Range("E10").Select() | Out-Null
